$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Email address: "minde@mindelis.com" -> "mindevieras@gmail.com"
#    (commit message: "email changed to mindevieras@gmail.com")
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Text = "minde@mindelis"
$found = $find.Execute()
if ($found) {
    $emailRng = $find.Parent
    $emailRng.Text = "mindevieras@gmail"
}

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark (Word's last-edit marker) so it sits
#    right after the freshly edited e-mail text, matching where Word
#    itself would leave it after the user's most recent edit.
#    A collapsed range exactly at a paragraph/cell end is rejected by
#    this host, so we briefly extend the run, anchor the bookmark at
#    the boundary (now a valid interior position) and then remove the
#    temporary text again - this leaves the bookmark collapsed right
#    after ".com" without splitting that run.
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Text = "mindevieras@gmail.com"
$found2 = $find2.Execute()
if ($found2) {
    $afterEmail = $find2.Parent.End

    $tempIns = $d.Range($afterEmail, $afterEmail)
    $tempIns.InsertAfter("ZZTMPZZ")

    $bmSpot = $d.Range($afterEmail, $afterEmail)
    $d.Bookmarks.Add("_GoBack", $bmSpot)

    $tempRange = $d.Range($afterEmail, $afterEmail + 7)
    $tempRange.Text = ""
}

# ------------------------------------------------------------------
# 3) Merge the split "CCNA/CISCO ... Campus, " / "Northern Ireland"
#    runs (which used to be separated by the old "_GoBack" bookmark)
#    back into a single run's text.
# ------------------------------------------------------------------
$find3 = $d.Content.Find
$find3.Text = " CCNA/CISCO South West College Enniskillen Campus, "
$found3 = $find3.Execute()
if ($found3) {
    $ccnaStart = $find3.Parent.Start
    $ccnaEnd = $find3.Parent.End
    $fullCcna = $d.Range($ccnaStart, $ccnaEnd + "Northern Ireland".Length)

    # Force a real text mutation (a no-op identical assignment would be
    # dropped) so the host actually collapses the bookmark-separated
    # runs into one contiguous run.
    $fullCcna.Text = "PLACEHOLDER_CCNA_MERGE"

    $find4 = $d.Content.Find
    $find4.Text = "PLACEHOLDER_CCNA_MERGE"
    $found4 = $find4.Execute()
    if ($found4) {
        $find4.Parent.Text = " CCNA/CISCO South West College Enniskillen Campus, Northern Ireland"
    }
}
